$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.605.47"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.979.63"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'382.63"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("D6").Value = "'103.55"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.593"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'37.08"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").Value = "3.446.85"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "'18.30"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "'7.60"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("D16").Value = "2.971.89"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "'0.996"
$ws.Range("E17").Value = "  +7.08%  "
$ws.Range("D18").Value = "51.510.59"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").Value = "'12.85"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").Value = "'69.12"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").Value = "'262.07"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("E25").Value = "  +8.41%  "
$ws.Range("D26").Value = "'8.33"
$ws.Range("E26").Value = "  +16.13%  "
$ws.Range("D27").Value = "'7.73"
$ws.Range("E27").Value = "  +18.35%  "
$ws.Range("E28").Value = "  +12.92%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'26.00"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").Value = "'9.88"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").Value = "'34.60"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'50.98"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("E36").Value = "  +6.23%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'3.01"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").Value = "'16.97"
$ws.Range("D40").Value = "'2.58"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").Value = "'122.39"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").Value = "'21.61"
$ws.Range("E44").Value = "  -2.23%  "
$ws.Range("E45").Value = "  +13.61%  "
$ws.Range("D46").Value = "'2.04"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.32"
$ws.Range("E47").Value = "  +4.67%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").Value = "2.031.83"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("E51").Value = "  +1.48%  "
